# "some adjustments based on risa's last comments" - figures adjusted for
# consistency: the %N data in columns E:F (rows 3-5) was bad, stats were
# rerun, and the new results were pasted back in without the "0.000"
# number format that the old (discarded) figures used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("foliar_cn")

# B3 already carries the plain (no custom number format) bordered style
# that the refreshed figures should use, so copy its format onto the
# E3:F5 block before writing the new values.
$ws.Range("B3").Copy()
$ws.Range("E3:F5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E3").Value = 0.889412488061744
$ws.Range("F3").Value = 0.349024916001268
$ws.Range("E4").Value = 4.2345550660915201
$ws.Range("F4").Value = 0.043502921115634703
$ws.Range("E5").Value = 0.59368054372753898
$ws.Range("F5").Value = 0.44370792261043501

# Match the refreshed selection left behind in the saved workbook.
[void]$ws.Range("K7").Select()
